$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct (explicit) new values per row, taken from the target diff.
# Columns updated: D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de comercializacion),
# R (Origen), S (Precio $/Kg), T (Kg / unidad)

$rowData = @{
    2  = @{ D = 44553; M = 200; N = 22000; O = 22000; P = 22000; Q = "`$/bandeja 6 kilos"; S = 3667; T = 6 }
    3  = @{ D = 44553; M = 150; N = 18000; O = 18000; P = 18000; Q = "`$/bandeja 6 kilos"; S = 3000; T = 6 }
    4  = @{ D = 44187; M = 45;  N = 14000; O = 14000; P = 14000; Q = "`$/bandeja 7 kilos"; S = 2000; T = 7 }
    5  = @{ D = 44187; M = 50;  N = 12000; O = 12000; P = 12000; Q = "`$/bandeja 7 kilos"; S = 1714; T = 7 }
    7  = @{ D = 44550; M = 60;  N = 24000; O = 24000; P = 24000; R = "Región Metropolitana"; S = 3429 }
    8  = @{ D = 44189; M = 20;  N = 15000; O = 15000; P = 15000; Q = "`$/bandeja 7 kilos"; S = 2143; T = 7 }
    9  = @{ D = 44189; M = 30;  N = 13000; O = 13000; P = 13000; Q = "`$/bandeja 7 kilos"; S = 1857; T = 7 }
    10 = @{ D = 44572; L = "Primera"; M = 65; N = 20000; O = 20000; P = 20000; Q = "`$/bandeja 6 kilos"; R = "Región Metropolitana"; S = 3333; T = 6 }
    11 = @{ D = 44204; M = 110; N = 7000; O = 7500; P = 7318; S = 1045 }
    12 = @{ D = 44558; L = "Especial"; M = 20; N = 22000; O = 22000; P = 22000; Q = "`$/bandeja 6 kilos"; R = "Provincia de San Felipe de Aconcagua"; S = 3667; T = 6 }
    13 = @{ D = 44558; M = 25;  N = 18000; O = 18000; P = 18000; R = "Provincia de San Felipe de Aconcagua"; S = 3000 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}
